$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from an existing row so the new date
# cells reuse the workbook's existing date style (numFmtId 14) instead
# of Excel synthesizing a brand-new custom number format.
$ws.Range("E129").Copy($ws.Range("E130"))
$ws.Range("E129").Copy($ws.Range("E131"))

# Row 130: ABC105 D, AC, failFirst=TRUE, date 2019-09-05, technique "二次元累積和"
$ws.Range("A130").Value = 105
$ws.Range("B130").Value = "D"
$ws.Range("C130").Value = "AC"
$ws.Range("D130").Value = $true
$ws.Range("E130").Value = 43713
$ws.Range("F130").Value = "二次元累積和"

# Row 131: ABC105 D, AC, failFirst=TRUE, date 2019-09-06, technique "平面走査、BIT"
$ws.Range("A131").Value = 105
$ws.Range("B131").Value = "D"
$ws.Range("C131").Value = "AC"
$ws.Range("D131").Value = $true
$ws.Range("E131").Value = 43714
$ws.Range("F131").Value = "平面走査、BIT"

$ws.Range("F132").Select()
